# Splits the single run of text in a paragraph into one run per
# whitespace-delimited word, with the separating spaces as their own
# runs (mirrors the OOXML diff: "foo bar" -> [w:r "foo"][w:r " "][w:r "bar"]).

function Split-ParagraphIntoWordRuns {
    param($Paragraph)

    $full = $Paragraph.Range
    # Exclude the trailing paragraph-mark character so the pPr (style) of
    # the host paragraph is left untouched by the XML replacement.
    $target = $word.ActiveDocument.Range($full.Start, $full.End - 1)

    $words = $target.Text -split " "

    $runsXml = New-Object System.Text.StringBuilder
    for ($i = 0; $i -lt $words.Length; $i++) {
        if ($i -gt 0) {
            [void]$runsXml.Append('<w:r><w:t xml:space="preserve"> </w:t></w:r>')
        }
        $escaped = $words[$i] -replace '&', '&amp;' -replace '<', '&lt;' -replace '>', '&gt;'
        [void]$runsXml.Append('<w:r><w:t xml:space="preserve">' + $escaped + '</w:t></w:r>')
    }

    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body><w:p>' + $runsXml.ToString() + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $target.InsertXML($xml)
}

$d = $word.ActiveDocument

$titlePara = $null
$authorPara = $null
$abstractPara = $null

foreach ($p in $d.Paragraphs) {
    $styleName = $p.Style.NameLocal
    if ($styleName -eq "Title" -and $titlePara -eq $null) {
        $titlePara = $p
    } elseif ($styleName -eq "Author" -and $authorPara -eq $null) {
        $authorPara = $p
    } elseif ($styleName -eq "Abstract" -and $abstractPara -eq $null) {
        $abstractPara = $p
    }
}

Split-ParagraphIntoWordRuns $titlePara
Split-ParagraphIntoWordRuns $authorPara
Split-ParagraphIntoWordRuns $abstractPara
